$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = '69+28=97'
$t.Cell(1, 2).Range.Text = '95-79=16'
$t.Cell(1, 3).Range.Text = '11-3=8'
$t.Cell(1, 4).Range.Text = '25-17=8'
$t.Cell(1, 5).Range.Text = '43-27=16'
$t.Cell(2, 1).Range.Text = '35+57=92'
$t.Cell(2, 2).Range.Text = '82-13=69'
$t.Cell(2, 3).Range.Text = '63-38=25'
$t.Cell(2, 4).Range.Text = '5+69=74'
$t.Cell(2, 5).Range.Text = '37+9=46'
$t.Cell(3, 1).Range.Text = '14+38=52'
$t.Cell(3, 2).Range.Text = '28+66=94'
$t.Cell(3, 3).Range.Text = '38+46=84'
$t.Cell(3, 4).Range.Text = '17+55=72'
$t.Cell(3, 5).Range.Text = '22+69=91'
$t.Cell(4, 1).Range.Text = '93-9=84'
$t.Cell(4, 2).Range.Text = '18+16=34'
$t.Cell(4, 3).Range.Text = '17+8=25'
$t.Cell(4, 4).Range.Text = '87-59=28'
$t.Cell(4, 5).Range.Text = '23-4=19'
$t.Cell(5, 1).Range.Text = '73+9=82'
$t.Cell(5, 2).Range.Text = '66-7=59'
$t.Cell(5, 3).Range.Text = '65+8=73'
$t.Cell(5, 4).Range.Text = '80-4=76'
$t.Cell(5, 5).Range.Text = '82-69=13'
$t.Cell(6, 1).Range.Text = '27+66=93'
$t.Cell(6, 2).Range.Text = '95-58=37'
$t.Cell(6, 3).Range.Text = '15+19=34'
$t.Cell(6, 4).Range.Text = '61-18=43'
$t.Cell(6, 5).Range.Text = '90-13=77'
$t.Cell(7, 1).Range.Text = '45+16=61'
$t.Cell(7, 2).Range.Text = '44+8=52'
$t.Cell(7, 3).Range.Text = '49+13=62'
$t.Cell(7, 4).Range.Text = '62-27=35'
$t.Cell(7, 5).Range.Text = '8+7=15'
$t.Cell(8, 1).Range.Text = '7+28=35'
$t.Cell(8, 2).Range.Text = '61-24=37'
$t.Cell(8, 3).Range.Text = '70-15=55'
$t.Cell(8, 4).Range.Text = '5+37=42'
$t.Cell(8, 5).Range.Text = '22-6=16'
$t.Cell(9, 1).Range.Text = '63-49=14'
$t.Cell(9, 2).Range.Text = '71-7=64'
$t.Cell(9, 3).Range.Text = '58+34=92'
$t.Cell(9, 4).Range.Text = '73-15=58'
$t.Cell(9, 5).Range.Text = '52-35=17'
$t.Cell(10, 1).Range.Text = '9+7=16'
$t.Cell(10, 2).Range.Text = '12-8=4'
$t.Cell(10, 3).Range.Text = '6+38=44'
$t.Cell(10, 4).Range.Text = '47+24=71'
$t.Cell(10, 5).Range.Text = '73-47=26'
$t.Cell(11, 1).Range.Text = '92-87=5'
$t.Cell(11, 2).Range.Text = '55-19=36'
$t.Cell(11, 3).Range.Text = '72-55=17'
$t.Cell(11, 4).Range.Text = '2+29=31'
$t.Cell(11, 5).Range.Text = '73-65=8'
$t.Cell(12, 1).Range.Text = '39+17=56'
$t.Cell(12, 2).Range.Text = '59+2=61'
$t.Cell(12, 3).Range.Text = '4+28=32'
$t.Cell(12, 4).Range.Text = '29+22=51'
$t.Cell(12, 5).Range.Text = '69+23=92'
$t.Cell(13, 1).Range.Text = '15+47=62'
$t.Cell(13, 2).Range.Text = '35-26=9'
$t.Cell(13, 3).Range.Text = '72-46=26'
$t.Cell(13, 4).Range.Text = '74-27=47'
$t.Cell(13, 5).Range.Text = '48+44=92'
$t.Cell(14, 1).Range.Text = '93-88=5'
$t.Cell(14, 2).Range.Text = '94-87=7'
$t.Cell(14, 3).Range.Text = '23+18=41'
$t.Cell(14, 4).Range.Text = '46+15=61'
$t.Cell(14, 5).Range.Text = '76+19=95'
$t.Cell(15, 1).Range.Text = '16+69=85'
$t.Cell(15, 2).Range.Text = '70-11=59'
$t.Cell(15, 3).Range.Text = '86-58=28'
$t.Cell(15, 4).Range.Text = '30-6=24'
$t.Cell(15, 5).Range.Text = '19+22=41'
$t.Cell(16, 1).Range.Text = '16+19=35'
$t.Cell(16, 2).Range.Text = '33-15=18'
$t.Cell(16, 3).Range.Text = '69+6=75'
$t.Cell(16, 4).Range.Text = '8+53=61'
$t.Cell(16, 5).Range.Text = '45-37=8'
$t.Cell(17, 1).Range.Text = '90-2=88'
$t.Cell(17, 2).Range.Text = '18+48=66'
$t.Cell(17, 3).Range.Text = '19+44=63'
$t.Cell(17, 4).Range.Text = '40-11=29'
$t.Cell(17, 5).Range.Text = '46+45=91'
$t.Cell(18, 1).Range.Text = '75-69=6'
$t.Cell(18, 2).Range.Text = '46+48=94'
$t.Cell(18, 3).Range.Text = '55-9=46'
$t.Cell(18, 4).Range.Text = '44+39=83'
$t.Cell(18, 5).Range.Text = '17+46=63'
$t.Cell(19, 1).Range.Text = '45+36=81'
$t.Cell(19, 2).Range.Text = '30-19=11'
$t.Cell(19, 3).Range.Text = '28+56=84'
$t.Cell(19, 4).Range.Text = '63-18=45'
$t.Cell(19, 5).Range.Text = '8+37=45'
$t.Cell(20, 1).Range.Text = '83-18=65'
$t.Cell(20, 2).Range.Text = '84-46=38'
$t.Cell(20, 3).Range.Text = '34-7=27'
$t.Cell(20, 4).Range.Text = '28-9=19'
$t.Cell(20, 5).Range.Text = '75+9=84'
